$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1351
$ws.Range("I103").Value = 2470
$ws.Range("J103").Value = 1071.25
$ws.Range("K103").Value = 7410
$ws.Range("L103").Value = 3213.75
$ws.Range("M103").Value = -6824
$ws.Range("N103").Value = -4385.75
$ws.Range("H112").Value = 5688.893
$ws.Range("J112").Value = 6283.56
$ws.Range("L112").Value = 18850.68
$ws.Range("N112").Value = -21066.68
$ws.Range("H113").Value = 9594.799999999999
$ws.Range("I113").Value = 29829.143
$ws.Range("J113").Value = 3436.5217
$ws.Range("K113").Value = 29829.143
$ws.Range("L113").Value = 3436.5217
$ws.Range("M113").Value = -26575.143
$ws.Range("N113").Value = -9944.521699999999
$ws.Range("H116").Value = 254740.16
$ws.Range("I116").Value = 3486.0715
$ws.Range("J116").Value = 390030.8
$ws.Range("K116").Value = 3486.0715
$ws.Range("L116").Value = 390030.8
$ws.Range("M116").Value = -44.07150000000001
$ws.Range("N116").Value = -396914.8
$ws.Range("H132").Value = 29253656
$ws.Range("I132").Value = 29863066
$ws.Range("K132").Value = 89589198
$ws.Range("M132").Value = -89586668
$ws.Range("H135").Value = 3540.2744
$ws.Range("I135").Value = 3826.3684
$ws.Range("J135").Value = 2704
$ws.Range("K135").Value = 34437.3156
$ws.Range("L135").Value = 24336
$ws.Range("M135").Value = -31902.3156
$ws.Range("N135").Value = -29406
$ws.Range("H138").Value = 1069.09
$ws.Range("I138").Value = 595.9706
$ws.Range("J138").Value = 2074.4688
$ws.Range("K138").Value = 1787.9118
$ws.Range("L138").Value = 6223.4064
$ws.Range("M138").Value = 3352.0882
$ws.Range("N138").Value = -16503.4064

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2318.64
$ws.Range("I32").Value = 2271.5483
$ws.Range("J32").Value = 2944.2856
$ws.Range("K32").Value = 2271.5483
$ws.Range("L32").Value = 2944.2856
$ws.Range("M32").Value = -1984.5483
$ws.Range("N32").Value = -3518.2856
$ws.Range("H61").Value = 8336943
$ws.Range("I61").Value = 11115378
$ws.Range("J61").Value = 1638.8
$ws.Range("K61").Value = 11115378
$ws.Range("L61").Value = 1638.8
$ws.Range("M61").Value = -11115166
$ws.Range("N61").Value = -2062.8
$ws.Range("H110").Value = 962.4
$ws.Range("I110").Value = 928.6667
$ws.Range("J110").Value = 1013
$ws.Range("K110").Value = 928.6667
$ws.Range("L110").Value = 1013
$ws.Range("M110").Value = 1116.3333
$ws.Range("N110").Value = -5103
$ws.Range("H132").Value = 3284394.5
$ws.Range("I132").Value = 5159960.5
$ws.Range("J132").Value = 2154.15
$ws.Range("K132").Value = 15479881.5
$ws.Range("L132").Value = 6462.450000000001
$ws.Range("M132").Value = -15477351.5
$ws.Range("N132").Value = -11522.45
$ws.Range("H136").Value = 8336943
$ws.Range("I136").Value = 11115378
$ws.Range("J136").Value = 1638.8
$ws.Range("K136").Value = 33346134
$ws.Range("L136").Value = 4916.4
$ws.Range("M136").Value = -33343584
$ws.Range("N136").Value = -10016.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4487.121
$ws.Range("I20").Value = 6207.2383
$ws.Range("K20").Value = 6207.2383
$ws.Range("M20").Value = -5960.2383
$ws.Range("H99").Value = 1045.5
$ws.Range("I99").Value = 1008.8333
$ws.Range("K99").Value = 1008.8333
$ws.Range("M99").Value = 489.1667
$ws.Range("H105").Value = 1723.6364
$ws.Range("I105").Value = 1723.6364
$ws.Range("K105").Value = 1723.6364
$ws.Range("M105").Value = 23.36359999999991
$ws.Range("H134").Value = 19623368
$ws.Range("I134").Value = 26687368
$ws.Range("J134").Value = 1144.4445
$ws.Range("K134").Value = 80062104
$ws.Range("L134").Value = 3433.3335
$ws.Range("M134").Value = -80059569
$ws.Range("N134").Value = -8503.333500000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2149363.2
$ws.Range("I58").Value = 2524405.8
$ws.Range("J58").Value = 11620.8
$ws.Range("K58").Value = 2524405.8
$ws.Range("L58").Value = 11620.8
$ws.Range("M58").Value = -2524202.8
$ws.Range("N58").Value = -12026.8
$ws.Range("H99").Value = 111113490
$ws.Range("I99").Value = 133335864
$ws.Range("J99").Value = 1607
$ws.Range("K99").Value = 133335864
$ws.Range("L99").Value = 1607
$ws.Range("M99").Value = -133334366
$ws.Range("N99").Value = -4603
$ws.Range("H126").Value = 111113490
$ws.Range("I126").Value = 133335864
$ws.Range("J126").Value = 1607
$ws.Range("K126").Value = 400007592
$ws.Range("L126").Value = 4821
$ws.Range("M126").Value = -400005122
$ws.Range("N126").Value = -9761
$ws.Range("H132").Value = 11500021
$ws.Range("I132").Value = 16669215
$ws.Range("J132").Value = 12924.444
$ws.Range("K132").Value = 50007645
$ws.Range("L132").Value = 38773.33199999999
$ws.Range("M132").Value = -50005115
$ws.Range("N132").Value = -43833.33199999999
$ws.Range("H134").Value = 13022817
$ws.Range("I134").Value = 13890711
$ws.Range("J134").Value = 10419133
$ws.Range("K134").Value = 41672133
$ws.Range("L134").Value = 31257399
$ws.Range("M134").Value = -41669598
$ws.Range("N134").Value = -31262469
$ws.Range("H136").Value = 2149363.2
$ws.Range("I136").Value = 2524405.8
$ws.Range("J136").Value = 11620.8
$ws.Range("K136").Value = 7573217.399999999
$ws.Range("L136").Value = 34862.39999999999
$ws.Range("M136").Value = -7570667.399999999
$ws.Range("N136").Value = -39962.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2096.1072
$ws.Range("I141").Value = 1529.2084
$ws.Range("J141").Value = 5497.5
$ws.Range("K141").Value = 4587.6252
$ws.Range("L141").Value = 16492.5
$ws.Range("M141").Value = 592.3747999999996
$ws.Range("N141").Value = -26852.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2833.3333
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2937.5
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2937.5
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4933.5
$ws.Range("H83").Value = 2833.3333
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2937.5
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 14687.5
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -24671.5
$ws.Range("H132").Value = 34725856
$ws.Range("I132").Value = 44445900
$ws.Range("J132").Value = 11404
$ws.Range("K132").Value = 133337700
$ws.Range("L132").Value = 34212
$ws.Range("M132").Value = -133335170
$ws.Range("N132").Value = -39272
$ws.Range("H134").Value = 16897.125
$ws.Range("J134").Value = 16897.125
$ws.Range("L134").Value = 50691.375
$ws.Range("N134").Value = -55761.375
$ws.Range("H135").Value = 35842.855
$ws.Range("J135").Value = 35842.855
$ws.Range("L135").Value = 35842.855
$ws.Range("N135").Value = -45982.855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4879755.5
$ws.Range("I132").Value = 9092223
$ws.Range("J132").Value = 2160.7368
$ws.Range("K132").Value = 27276669
$ws.Range("L132").Value = 6482.2104
$ws.Range("M132").Value = -27274139
$ws.Range("N132").Value = -11542.2104
$ws.Range("H136").Value = 17547016
$ws.Range("I136").Value = 22225698
$ws.Range("J136").Value = 1958.75
$ws.Range("K136").Value = 66677094
$ws.Range("L136").Value = 5876.25
$ws.Range("M136").Value = -66674544
$ws.Range("N136").Value = -10976.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 37902530
$ws.Range("I136").Value = 33102686
$ws.Range("J136").Value = 62501710
$ws.Range("K136").Value = 99308058
$ws.Range("L136").Value = 187505130
$ws.Range("M136").Value = -99305508
$ws.Range("N136").Value = -187510230
